# Weekly data refresh: insert a new week's worth of Coliflor price data
# (one "Primera" row and one "Segunda" row) right before the existing
# row that used to be row 257, pushing all the following rows down by
# two positions. This matches the pattern already used throughout the
# sheet where each week contributes exactly two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 257-258; everything from the old row 257
# onward shifts down to 259 onward. Excel copies the formatting
# (including the date number format on column D) from the row above,
# just like a normal "Insert Copied Cells"/"Insert Sheet Rows" action.
$ws.Range("A257:A258").EntireRow.Insert()

# --- New row 257 ("Primera") ---
$ws.Range("A257").Value = 8
$ws.Range("B257").Value = "Terminal La Palmera de La Serena"
$ws.Range("C257").Value = "Coquimbo"
$ws.Range("D257").Value = 44460
$ws.Range("E257").Value = 4
$ws.Range("F257").Value = 100112008
$ws.Range("G257").Value = "Coliflor"
$ws.Range("H257").Value = "Sin especificar"
$ws.Range("I257").Value = "Primera"
$ws.Range("J257").Value = 2200
$ws.Range("K257").Value = 600
$ws.Range("L257").Value = 700
$ws.Range("M257").Value = 650
$ws.Range("N257").Value = "`$/unidad"
$ws.Range("O257").Value = "Provincia del Elquí"
$ws.Range("P257").Value = 650
$ws.Range("Q257").Value = 1
$ws.Range("R257").Value = "Hortaliza"

# --- New row 258 ("Segunda") ---
$ws.Range("A258").Value = 8
$ws.Range("B258").Value = "Terminal La Palmera de La Serena"
$ws.Range("C258").Value = "Coquimbo"
$ws.Range("D258").Value = 44460
$ws.Range("E258").Value = 4
$ws.Range("F258").Value = 100112008
$ws.Range("G258").Value = "Coliflor"
$ws.Range("H258").Value = "Sin especificar"
$ws.Range("I258").Value = "Segunda"
$ws.Range("J258").Value = 1400
$ws.Range("K258").Value = 500
$ws.Range("L258").Value = 550
$ws.Range("M258").Value = 525
$ws.Range("N258").Value = "`$/unidad"
$ws.Range("O258").Value = "Provincia del Elquí"
$ws.Range("P258").Value = 525
$ws.Range("Q258").Value = 1
$ws.Range("R258").Value = "Hortaliza"
